# GQA - Avaliação (Quase Pronta)
# Adds new "Projeto 1" (column D) evidence markers / hyperlinks mirroring the
# existing "Processo" (column C) ones, plus two brand-new evidence rows
# (61/62 under RAP 6, mirrored partially at 75 under RAP 8, and 82 under
# RAP 10) that reference additional supporting documents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GQA")

# --- New rows: fill in column A labels (new shared strings) --------------
# Order matters only for shared-string insertion order, but doesn't affect
# correctness either way.
$ws.Range("A82").Value = "Certificados"
$ws.Range("A61").Value = "GQA - Especificação da Garantia da Qualidade (Tabela de Reuniões)"
$ws.Range("A75").Value = "GQA - Especificação da Garantia da Qualidade (Tabela de Reuniões)"
$ws.Range("A62").Value = "GQA - Documento de Integrantes"

# --- Column D: mark "X" for rows that previously had no project evidence -
$ws.Range("D61").Value = "X"
$ws.Range("D62").Value = "X"
$ws.Range("D75").Value = "X"
$ws.Range("D82").Value = "X"

# --- Hyperlinks on column D (new + pre-existing "X" cells) ---------------
# Existing cells D6, D13, D20, D28, D54, D67, D81 already show "X"; they
# simply gain a hyperlink (which also switches their style to "Hiperlink").
$baseTemplates = "https://github.com/samuelrcosta/processo-software/blob/master/Raiz/Processo/Templates/"
$baseDefinicao = "https://github.com/samuelrcosta/processo-software/blob/master/Raiz/Processo/Defini%C3%A7%C3%A3o/GQA-Processo.docx"
$baseCertificados = "https://github.com/samuelrcosta/processo-software/tree/master/Raiz/Projeto/Garantia%20da%20Qualidade/GQA%20-%20Certificados"
$especificacao = $baseTemplates + "GQA%20-%20Especifica%C3%A7%C3%A3o%20da%20Garantia%20da%20Qualidade.docx"
$integrantes = $baseTemplates + "GQA%20-%20Documento%20de%20Integrantes.docx"
$naoConformidades = $baseTemplates + "GQA%20-%20Documento%20de%20N%C3%A3o%20Conformidades.docx"

$ws.Hyperlinks.Add($ws.Range("D54"), $integrantes)
$ws.Hyperlinks.Add($ws.Range("D67"), $integrantes)
$ws.Hyperlinks.Add($ws.Range("D81"), $especificacao)
$ws.Hyperlinks.Add($ws.Range("D13"), $especificacao)
$ws.Hyperlinks.Add($ws.Range("D6"), $especificacao)
$ws.Hyperlinks.Add($ws.Range("D82"), $baseCertificados)
$ws.Hyperlinks.Add($ws.Range("D75"), $especificacao)
$ws.Hyperlinks.Add($ws.Range("D61"), $especificacao)
$ws.Hyperlinks.Add($ws.Range("D62"), $integrantes)
$ws.Hyperlinks.Add($ws.Range("D28"), $naoConformidades)
$ws.Hyperlinks.Add($ws.Range("D20"), $baseDefinicao)

# --- View state: active cell ends up on D20 -------------------------------
$ws.Range("D20").Select()
